$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 2020 value in L3 to 2021 (year header correction)
$ws.Range("L3").Value = 2021

# Add a new column M: copy formatting from the adjacent column L, then set values
# Row 3 header (year 2022)
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("M3").Value = 2022

# Row 4 data value (6.18)
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 6.18

# Update the current selection to match the new active cell (M9)
$ws.Range("M9").Select()
